$d = $word.ActiveDocument

# The schedule line reads "2021 04 13 1300-1430" (date 04/13) and the
# session date moved to the 14th, i.e. "2021 04 14 1300-1430".
# "13 1300-1430" is unique in the document, so target it precisely and
# only swap the day number, leaving the rest of the text untouched.
$d.Content.Find.Execute("13 1300-1430", $true, $false, $false, $false, $false,
                         $true, 1, $false, "14 1300-1430", 2)
